# DB update for financials
# Adds a second "PERIOD/YEAR/CLUSTER/ACCOUNT NAME/LOCATION" block (columns M,O,Q,S,U)
# next to the existing Financials table, plus a duplicated header+data block
# further down the sheet (rows 24, 26-34) for the 2021-2022 period.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Financials")
$ws3 = $wb.Worksheets.Item("Fulfilment_view")

# --- Write new cell values in the exact order the original author typed them,
#     so brand-new shared-string entries land at the same indices as the
#     target workbook (25 "2020-2021" .. 37 "Offsite"). ---

$ws1.Range("O8").Value  = "2020-2021"
$ws1.Range("U8").Value  = "Onsite"
$ws1.Range("Q8").Value  = "UKISA"
$ws1.Range("M6").Value  = "PERIOD"
$ws1.Range("O6").Value  = "YEAR"
$ws1.Range("Q6").Value  = "CLUSTER"
$ws1.Range("S6").Value  = "ACCOUNT NAME"
$ws1.Range("U6").Value  = "LOCATION"
$ws1.Range("S8").Value  = "Finastra"
$ws1.Range("O26").Value = "2021-2022"
$ws1.Range("Q26").Value = "ERGER"
$ws1.Range("S26").Value = "GER"
$ws1.Range("U26").Value = "Offsite"

# --- Fill the rest of the first block (rows 8-16): PERIOD/YEAR = 2020,
#     CLUSTER = UKISA, ACCOUNT NAME = Finastra, LOCATION = Onsite. ---

$ws1.Range("M8:M16").Value = 2020
$ws1.Range("O9:O16").Value = "2020-2021"
$ws1.Range("Q9:Q16").Value = "UKISA"
$ws1.Range("S9:S16").Value = "Finastra"
$ws1.Range("U9:U16").Value = "Onsite"

# --- Repeat the header row further down the sheet (row 24) ---

$ws1.Range("C24").Value = "id"
$ws1.Range("E24").Value = "label"
$ws1.Range("G24").Value = "percentage"
$ws1.Range("I24").Value = "title"
$ws1.Range("K24").Value = "heading"
$ws1.Range("M24").Value = "PERIOD"
$ws1.Range("O24").Value = "YEAR"
$ws1.Range("Q24").Value = "CLUSTER"
$ws1.Range("S24").Value = "ACCOUNT NAME"
$ws1.Range("U24").Value = "LOCATION"

# --- Second data block (rows 26-34), mirroring rows 8-16 for 2021-2022 ---

$ws1.Range("C26").Value = 1
$ws1.Range("E26").Value = "Budget"
$ws1.Range("G26").Value = 55
$ws1.Range("I26").Value = "Revenue"
$ws1.Range("K26").Value = "Total Revenue: 2023-2024"

$ws1.Range("C27").Value = 2
$ws1.Range("E27").Value = "Actuals"
$ws1.Range("G27").Value = 4
$ws1.Range("I27").Value = "Revenue"
$ws1.Range("K27").Value = "Total Revenue: 2023-2024"

$ws1.Range("C28").Value = 3
$ws1.Range("E28").Value = "Achieved percentage"
$ws1.Range("G28").Value = 67
$ws1.Range("I28").Value = "Revenue"
$ws1.Range("K28").Value = "Total Revenue: 2023-2024"

$ws1.Range("C29").Value = 4
$ws1.Range("E29").Value = "Budget"
$ws1.Range("G29").Value = 8
$ws1.Range("I29").Value = "Gross Margin $"
$ws1.Range("K29").Value = "Total Revenue: 2023-2024"

$ws1.Range("C30").Value = 5
$ws1.Range("E30").Value = "Actuals"
$ws1.Range("G30").Value = 99
$ws1.Range("I30").Value = "Gross Margin $"
$ws1.Range("K30").Value = "Total Revenue: 2023-2024"

$ws1.Range("C31").Value = 6
$ws1.Range("E31").Value = "Achieved percentage"
$ws1.Range("G31").Value = 7
$ws1.Range("I31").Value = "Gross Margin $"
$ws1.Range("K31").Value = "Total Revenue: 2023-2024"

$ws1.Range("C32").Value = 7
$ws1.Range("E32").Value = "Budget"
$ws1.Range("G32").Value = 56
$ws1.Range("I32").Value = "Gross Margin %"
$ws1.Range("K32").Value = "Total Revenue: 2023-2024"

$ws1.Range("C33").Value = 8
$ws1.Range("E33").Value = "Actuals"
$ws1.Range("G33").Value = 32
$ws1.Range("I33").Value = "Gross Margin %"
$ws1.Range("K33").Value = "Total Revenue: 2023-2024"

$ws1.Range("C34").Value = 9
$ws1.Range("E34").Value = "Achieved percentage"
$ws1.Range("G34").Value = 67
$ws1.Range("I34").Value = "Gross Margin %"
$ws1.Range("K34").Value = "Total Revenue: 2023-2024"

$ws1.Range("M26:M34").Value = 2021
$ws1.Range("O27:O34").Value = "2021-2022"
$ws1.Range("Q27:Q34").Value = "ERGER"
$ws1.Range("S27:S34").Value = "GER"
$ws1.Range("U27:U34").Value = "Offsite"

# --- Column widths for the new columns (closest achievable values) ---

$ws1.Columns("E").ColumnWidth = 17.5
$ws1.Columns("K").ColumnWidth = 21.5
$ws1.Columns("O").ColumnWidth = 9.5

# --- View / selection state ---
# (select on Fulfilment_view first, then activate+select Financials last, so
#  Financials ends up as the sheet that is actually active/tabSelected.)

$ws3.Range("P7").Select()

$ws1.Activate()
$ws1.Range("N25").Select()
